$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D, shifting the existing D:K quarterly
# data to E:L (making room for a newer quarter of figures in column D).
$ws.Columns("D:D").Insert()

# The freshly inserted column has default formatting; clone the formatting
# (number formats / styles) from the now-adjacent column E so the new D
# column matches the rest of its row (date format on the "Period Ending"
# rows, the thousands number format elsewhere, etc.). Restrict the copy to
# the three data blocks that actually carried a column D cell before the
# insert (the blank section-header rows 5, 6, 37 and 79 never had one).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

# New column D values (most recent quarter) for each data row.
$ws.Range("D7").Value = 43373
$ws.Range("D8").Value = 7500
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = -200
$ws.Range("D17").Value = 700
$ws.Range("D18").Value = 6800
$ws.Range("D20").Value = -4700
$ws.Range("D21").Value = 2700
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 2100
$ws.Range("D24").Value = 400
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 1700
$ws.Range("D27").Value = 1700
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 4700
$ws.Range("D33").Value = 1700
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 1700

$ws.Range("D38").Value = 43373
$ws.Range("D41").Value = 8600
$ws.Range("D42").Value = 30700
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 20700
$ws.Range("D49").Value = 7900
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 2600
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 679700
$ws.Range("D57").Value = 200
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 605200
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 35100
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 74400
$ws.Range("D77").Value = 0

$ws.Range("D80").Value = 43373
$ws.Range("D81").Value = 1700
$ws.Range("D83").Value = 500
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 1900
$ws.Range("D91").Value = -1800
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 35100
$ws.Range("D96").Value = -600
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -14600
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 22500

# Row 15 also picked up revised "NA" markers (shifted from the old D:I=0
# values) across E:J, leaving only the last two quarters (K15/L15) as 0.
$ws.Range("E15:J15").Value = "NA"
